$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Background:" Instructions paragraph -> "Interpretation:" paragraph
#    - paragraph-mark run properties change from <i/><iCs/> to
#      <b val="0"/><bCs val="0"/><i/><iCs/>
#    - runs/text change entirely (keeps the lastRenderedPageBreak hint
#      on the first run)
# ---------------------------------------------------------------------
$introPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.StartsWith("Background:")) {
        $introPara = $pp
        break
    }
}

$introFull = $introPara.Range
$introXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="006949B6" w:rsidRPr="00D420BF" w:rsidRDefault="006949B6" w:rsidP="00D420BF"><w:pPr><w:pStyle w:val="Instructions"/><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Interpretation: </w:t></w:r><w:r><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:i/><w:iCs/></w:rPr><w:t>Read through the passage, and then review and make notes under the questions below.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$introFull.InsertXML($introXml)

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from inside the "List Paul's prayer
#    requests..." question (where it split two runs) to the very start
#    of the "Note the one thing..." question paragraph.
# ---------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$notePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.StartsWith("Not") -and $pp.Style.NameLocal -eq "Question") {
        $notePara = $pp
        break
    }
}
$noteStart = $notePara.Range.Start
$bmRange = $d.Range($noteStart, $noteStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 3) The "List Paul's prayer requests..." question used to be split
#    into two runs by the bookmark; now that the bookmark has moved
#    away, re-write it as a single contiguous run with the same text.
# ---------------------------------------------------------------------
$prayerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.StartsWith("List Paul")) {
        $prayerPara = $pp
        break
    }
}
$prayerFull = $prayerPara.Range
$prayerText = $prayerFull.Text.TrimEnd([char]13, [char]7)
$delRange = $d.Range($prayerFull.Start, $prayerFull.End)
$delRange.Text = ""
$insPoint = $d.Range($prayerFull.Start, $prayerFull.Start)
$insPoint.InsertAfter($prayerText)
